# Apply the commit's changes:
#  - Clear A1 (was "#")
#  - Clear the "Maßstab"/"Gym" column (C) for the subject-scale block
#    (rows 8, 9, 11-15) and merge B:C for those rows so the label in B
#    spans the now-empty C column.
#  - Move the active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray "#" header cell.
$ws.Range("A1").ClearContents()

# Clear column C in the affected rows, then merge B:C so the row label
# (Schüler(in)/id scale, student names, …) spans across.
$rows = @(8, 9, 11, 12, 13, 14, 15)
foreach ($r in $rows) {
    $ws.Range("C$r").ClearContents()
    $ws.Range("B$r`:C$r").Merge()
}

# Reset the saved selection to A1.
$ws.Range("A1").Select()
